$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("H132").Value = 9529570
$ws.Range("I132").Value = 12827294
$ws.Range("J132").Value = 2810.7778
$ws.Range("K132").Value = 38481882
$ws.Range("L132").Value = 8432.3334
$ws.Range("M132").Value = -38479352
$ws.Range("N132").Value = -13492.3334
$ws.Range("H138").Value = 544830.2
$ws.Range("I138").Value = 975.87177
$ws.Range("J138").Value = 1062155.1
$ws.Range("K138").Value = 2927.61531
$ws.Range("L138").Value = 3186465.3
$ws.Range("M138").Value = 2212.38469
$ws.Range("N138").Value = -3196745.3
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3480.3015
$ws.Range("I32").Value = 3205.353
$ws.Range("J32").Value = 4648.8335
$ws.Range("K32").Value = 3205.353
$ws.Range("L32").Value = 4648.8335
$ws.Range("M32").Value = -2918.353
$ws.Range("N32").Value = -5222.8335
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30676
$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32340
$ws.Range("H96").Value = 7672
$ws.Range("J96").Value = 7672
$ws.Range("L96").Value = 7672
$ws.Range("N96").Value = -13164
$ws.Range("H97").Value = 454.33334
$ws.Range("I97").Value = 401.07144
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 401.07144
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = 94.92856
$ws.Range("N97").Value = -2192
$ws.Range("H132").Value = 2250.875
$ws.Range("I132").Value = 2077.077
$ws.Range("J132").Value = 2456.2727
$ws.Range("K132").Value = 6231.231000000001
$ws.Range("L132").Value = 7368.8181
$ws.Range("M132").Value = -3701.231000000001
$ws.Range("N132").Value = -12428.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 22728262
$ws.Range("I94").Value = 31250812
$ws.Range("K94").Value = 31250812
$ws.Range("M94").Value = -31250361
$ws.Range("H105").Value = 111114250
$ws.Range("I105").Value = 142860030
$ws.Range("J105").Value = 3999.5
$ws.Range("K105").Value = 142860030
$ws.Range("L105").Value = 3999.5
$ws.Range("M105").Value = -142858283
$ws.Range("N105").Value = -7493.5
$ws.Range("H134").Value = 10032.143
$ws.Range("I134").Value = 1731.375
$ws.Range("J134").Value = 21099.834
$ws.Range("K134").Value = 5194.125
$ws.Range("L134").Value = 63299.50199999999
$ws.Range("M134").Value = -2659.125
$ws.Range("N134").Value = -68369.50199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 40001144
$ws.Range("I16").Value = 52632590
$ws.Range("J16").Value = 1552.1666
$ws.Range("K16").Value = 52632590
$ws.Range("L16").Value = 1552.1666
$ws.Range("M16").Value = -52632303
$ws.Range("N16").Value = -2126.1666
$ws.Range("H31").Value = 1719.2
$ws.Range("I31").Value = 1345.6
$ws.Range("J31").Value = 2092.8
$ws.Range("K31").Value = 1345.6
$ws.Range("L31").Value = 2092.8
$ws.Range("M31").Value = -1050.6
$ws.Range("N31").Value = -2682.8
$ws.Range("H34").Value = 1719.2
$ws.Range("I34").Value = 1345.6
$ws.Range("J34").Value = 2092.8
$ws.Range("K34").Value = 1345.6
$ws.Range("L34").Value = 2092.8
$ws.Range("M34").Value = -1143.6
$ws.Range("N34").Value = -2496.8
$ws.Range("H105").Value = 722.3889
$ws.Range("I105").Value = 710.2308
$ws.Range("J105").Value = 754
$ws.Range("K105").Value = 710.2308
$ws.Range("L105").Value = 754
$ws.Range("M105").Value = 1036.7692
$ws.Range("N105").Value = -4248
$ws.Range("H111").Value = 47499.5
$ws.Range("J111").Value = 47499.5
$ws.Range("L111").Value = 47499.5
$ws.Range("N111").Value = -55679.5
$ws.Range("H113").Value = 40001144
$ws.Range("I113").Value = 52632590
$ws.Range("J113").Value = 1552.1666
$ws.Range("K113").Value = 52632590
$ws.Range("L113").Value = 1552.1666
$ws.Range("M113").Value = -52630420
$ws.Range("N113").Value = -5892.1666
$ws.Range("H132").Value = 7288.45
$ws.Range("I132").Value = 9826.75
$ws.Range("K132").Value = 29480.25
$ws.Range("M132").Value = -26950.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 314
$ws.Range("I13").Value = 179.6
$ws.Range("J13").Value = 650
$ws.Range("K13").Value = 538.8
$ws.Range("L13").Value = 1950
$ws.Range("M13").Value = -370.8
$ws.Range("N13").Value = -2286
$ws.Range("H107").Value = 8078.6924
$ws.Range("I107").Value = 296.8
$ws.Range("J107").Value = 12942.375
$ws.Range("K107").Value = 890.4000000000001
$ws.Range("L107").Value = 38827.125
$ws.Range("M107").Value = 1029.6
$ws.Range("N107").Value = -42667.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3850
$ws.Range("I80").Value = 1800
$ws.Range("J80").Value = 5900
$ws.Range("K80").Value = 1800
$ws.Range("L80").Value = 5900
$ws.Range("M80").Value = -802
$ws.Range("N80").Value = -7896
$ws.Range("H83").Value = 3850
$ws.Range("I83").Value = 1800
$ws.Range("J83").Value = 5900
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 29500
$ws.Range("M83").Value = -4008
$ws.Range("N83").Value = -39484
$ws.Range("H97").Value = 835
$ws.Range("I97").Value = 820
$ws.Range("K97").Value = 820
$ws.Range("M97").Value = -324
$ws.Range("H126").Value = 2113.6365
$ws.Range("I126").Value = 1778.5714
$ws.Range("J126").Value = 2700
$ws.Range("K126").Value = 5335.7142
$ws.Range("L126").Value = 8100
$ws.Range("M126").Value = -2865.7142
$ws.Range("N126").Value = -13040
$ws.Range("H132").Value = 2655.5518
$ws.Range("I132").Value = 2248.4119
$ws.Range("J132").Value = 3232.3333
$ws.Range("K132").Value = 6745.2357
$ws.Range("L132").Value = 9696.999899999999
$ws.Range("M132").Value = -4215.2357
$ws.Range("N132").Value = -14756.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1580.3334
$ws.Range("J22").Value = 1875
$ws.Range("L22").Value = 1875
$ws.Range("N22").Value = -2465
$ws.Range("H27").Value = 1580.3334
$ws.Range("J27").Value = 1875
$ws.Range("L27").Value = 1875
$ws.Range("N27").Value = -2089
$ws.Range("H68").Value = 2065.8
$ws.Range("I68").Value = 1689.8889
$ws.Range("J68").Value = 2629.6667
$ws.Range("K68").Value = 1689.8889
$ws.Range("L68").Value = 2629.6667
$ws.Range("M68").Value = -940.8888999999999
$ws.Range("N68").Value = -4127.6667
$ws.Range("H71").Value = 2065.8
$ws.Range("I71").Value = 1689.8889
$ws.Range("J71").Value = 2629.6667
$ws.Range("K71").Value = 8449.4445
$ws.Range("L71").Value = 13148.3335
$ws.Range("M71").Value = -4705.4445
$ws.Range("N71").Value = -20636.3335
$ws.Range("H82").Value = 2066.5
$ws.Range("I82").Value = 2125
$ws.Range("K82").Value = 2125
$ws.Range("M82").Value = -1764
$ws.Range("H85").Value = 2066.5
$ws.Range("I85").Value = 2125
$ws.Range("K85").Value = 2125
$ws.Range("M85").Value = -877
$ws.Range("H104").Value = 5458.2856
$ws.Range("J104").Value = 5458.2856
$ws.Range("L104").Value = 5458.2856
$ws.Range("N104").Value = -12446.2856
$ws.Range("H132").Value = 31582.53
$ws.Range("I132").Value = 1753.5333
$ws.Range("J132").Value = 55131.74
$ws.Range("K132").Value = 5260.5999
$ws.Range("L132").Value = 165395.22
$ws.Range("M132").Value = -2730.5999
$ws.Range("N132").Value = -170455.22
$ws.Range("H136").Value = 5211.2964
$ws.Range("I136").Value = 5709.1304
$ws.Range("K136").Value = 17127.3912
$ws.Range("M136").Value = -14577.3912

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2000.6666
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 2000.6666
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 4001.3332
$ws.Range("N81").Value = -6123.3332
$ws.Range("H84").Value = 2000.6666
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 2000.6666
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 20006.666
$ws.Range("N84").Value = -30614.666
$ws.Range("H122").Value = 20001074
$ws.Range("I122").Value = 20001074
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 60003222
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -60000772
$ws.Range("H136").Value = 624.0909
$ws.Range("I136").Value = 432.94736
$ws.Range("K136").Value = 1298.84208
$ws.Range("M136").Value = 1251.15792
$ws.Range("M81").ClearContents()
$ws.Range("M84").ClearContents()
$ws.Range("N122").ClearContents()

Write-Host "Applied all Kujata_Profits market data updates"